$wb = $excel.ActiveWorkbook

# The update touches both the "展览" sheet and the "全部类型" sheet,
# which mirror the same data rows.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 9
    $ws.Range("F5").Value = 71
    $ws.Range("F6").Value = 15
}
